$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A100").Value = 15
$ws.Range("B100").Value = 4
$ws.Range("C100").Value = 4

$ws.Range("A101").Value = 15
$ws.Range("B101").Value = 7
$ws.Range("C101").Value = 5
